$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Anas platyrhynchos"
$ws.Range("H2").Value = "Domestic duck"

$ws.Range("G35").Value = "Spermophilus beecheyi"
$ws.Range("H35").Value = "California ground squirrel"
